$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 677.1875
$ws.Range("I15").Value = 677.1875
$ws.Range("K15").Value = 2031.5625
$ws.Range("M15").Value = -1862.5625
$ws.Range("H33").Value = 142.81818
$ws.Range("I33").Value = 134.15
$ws.Range("K33").Value = 134.15
$ws.Range("M33").Value = 94.84999999999999
$ws.Range("H62").Value = 6440.857
$ws.Range("J62").Value = 6558
$ws.Range("L62").Value = 6558
$ws.Range("N62").Value = -7806
$ws.Range("H65").Value = 6440.857
$ws.Range("J65").Value = 6558
$ws.Range("L65").Value = 32790
$ws.Range("N65").Value = -39030
$ws.Range("H70").Value = 3738.3076
$ws.Range("I70").Value = 3582.6667
$ws.Range("K70").Value = 10748.0001
$ws.Range("M70").Value = -10478.0001
$ws.Range("H73").Value = 3738.3076
$ws.Range("I73").Value = 3582.6667
$ws.Range("K73").Value = 10748.0001
$ws.Range("M73").Value = -9812.000100000001
$ws.Range("H92").Value = 443.4375
$ws.Range("I92").Value = 413
$ws.Range("K92").Value = 413
$ws.Range("M92").Value = 835
$ws.Range("H106").Value = 5000
$ws.Range("I106").Value = 5000
$ws.Range("K106").Value = 5000
$ws.Range("M106").Value = -4369
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("H137").Value = 4041.25
$ws.Range("J137").Value = 4633
$ws.Range("L137").Value = 13899
$ws.Range("N137").Value = -18999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 500
$ws.Range("I4").Value = 500
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 500
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -384
$ws.Range("H5").Value = 232.83333
$ws.Range("I5").Value = 232.83333
$ws.Range("K5").Value = 232.83333
$ws.Range("M5").Value = -120.83333
$ws.Range("H32").Value = 1964.1471
$ws.Range("I32").Value = 1720.6364
$ws.Range("K32").Value = 1720.6364
$ws.Range("M32").Value = -1433.6364
$ws.Range("H63").Value = 2095.5
$ws.Range("I63").Value = 2154.6
$ws.Range("J63").Value = 1800
$ws.Range("K63").Value = 2154.6
$ws.Range("L63").Value = 1800
$ws.Range("M63").Value = -1468.6
$ws.Range("N63").Value = -3172
$ws.Range("H66").Value = 2095.5
$ws.Range("I66").Value = 2154.6
$ws.Range("J66").Value = 1800
$ws.Range("K66").Value = 10773
$ws.Range("L66").Value = 9000
$ws.Range("M66").Value = -7341
$ws.Range("N66").Value = -15864
$ws.Range("H74").Value = 3845.2
$ws.Range("I74").Value = 4670.6665
$ws.Range("J74").Value = 2607
$ws.Range("K74").Value = 4670.6665
$ws.Range("L74").Value = 2607
$ws.Range("M74").Value = -3796.6665
$ws.Range("N74").Value = -4355
$ws.Range("H77").Value = 3845.2
$ws.Range("I77").Value = 4670.6665
$ws.Range("J77").Value = 2607
$ws.Range("K77").Value = 23353.3325
$ws.Range("L77").Value = 13035
$ws.Range("M77").Value = -18985.3325
$ws.Range("N77").Value = -21771
$ws.Range("H122").Value = 2800
$ws.Range("I122").Value = 2850
$ws.Range("J122").Value = 2750
$ws.Range("K122").Value = 8550
$ws.Range("L122").Value = 8250
$ws.Range("M122").Value = -6100
$ws.Range("N122").Value = -13150
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").ClearContents()
$ws.Range("N139").Value = 0

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 232.83333
$ws.Range("I4").Value = 232.83333
$ws.Range("K4").Value = 232.83333
$ws.Range("M4").Value = -117.83333
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H86").Value = 3647.7144
$ws.Range("I86").Value = 5450
$ws.Range("K86").Value = 5450
$ws.Range("M86").Value = -4327
$ws.Range("H89").Value = 3647.7144
$ws.Range("I89").Value = 5450
$ws.Range("K89").Value = 27250
$ws.Range("M89").Value = -21634
$ws.Range("H99").Value = 3178.25
$ws.Range("I99").Value = 3103.5454
$ws.Range("K99").Value = 3103.5454
$ws.Range("M99").Value = -1605.5454
$ws.Range("H105").Value = 1743.4762
$ws.Range("I105").Value = 1715.65
$ws.Range("K105").Value = 1715.65
$ws.Range("M105").Value = 31.34999999999991
$ws.Range("H134").Value = 4215.08
$ws.Range("I134").Value = 4215.08
$ws.Range("K134").Value = 12645.24
$ws.Range("M134").Value = -10110.24

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 7540.625
$ws.Range("I25").Value = 5049.8335
$ws.Range("K25").Value = 5049.8335
$ws.Range("M25").Value = -4875.8335
$ws.Range("H31").Value = 2758.6667
$ws.Range("I31").Value = 2758.6667
$ws.Range("K31").Value = 2758.6667
$ws.Range("M31").Value = -2463.6667
$ws.Range("H34").Value = 2758.6667
$ws.Range("I34").Value = 2758.6667
$ws.Range("K34").Value = 2758.6667
$ws.Range("M34").Value = -2556.6667
$ws.Range("H50").Value = 20178.834
$ws.Range("J50").Value = 19998
$ws.Range("L50").Value = 19998
$ws.Range("N50").Value = -21248
$ws.Range("H105").Value = 925.125
$ws.Range("I105").Value = 967.6667
$ws.Range("K105").Value = 967.6667
$ws.Range("M105").Value = 779.3333

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 493.16666
$ws.Range("I117").Value = 549
$ws.Range("J117").Value = 465.25
$ws.Range("K117").Value = 1647
$ws.Range("L117").Value = 1395.75
$ws.Range("M117").Value = 1795
$ws.Range("N117").Value = -8279.75
$ws.Range("H122").Value = 1132.6666
$ws.Range("I122").Value = 1347.75
$ws.Range("J122").Value = 960.6
$ws.Range("K122").Value = 12129.75
$ws.Range("L122").Value = 8645.4
$ws.Range("M122").Value = -9679.75
$ws.Range("N122").Value = -13545.4
$ws.Range("H132").Value = 3012.5
$ws.Range("J132").Value = 3850
$ws.Range("L132").Value = 34650
$ws.Range("N132").Value = -39710
$ws.Range("H139").Value = 1250.7142
$ws.Range("I139").Value = 1250.7142
$ws.Range("K139").Value = 3752.1426
$ws.Range("M139").Value = 1387.8574

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 1200044.9
$ws.Range("I11").Value = 257213.28
$ws.Range("J11").Value = 2850000
$ws.Range("K11").Value = 257213.28
$ws.Range("L11").Value = 2850000
$ws.Range("M11").Value = -257074.28
$ws.Range("N11").Value = -2850278
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("H122").Value = 880
$ws.Range("I122").Value = 880
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2640
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -190
$ws.Range("H126").Value = 6479.8
$ws.Range("I126").Value = 7599.75
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 22799.25
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -20329.25
$ws.Range("N126").Value = -10940

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 25668.334
$ws.Range("J11").Value = 25007
$ws.Range("L11").Value = 25007
$ws.Range("N11").Value = -25287
$ws.Range("H22").Value = 771.2857
$ws.Range("I22").Value = 670.75
$ws.Range("K22").Value = 670.75
$ws.Range("M22").Value = -375.75
$ws.Range("H27").Value = 771.2857
$ws.Range("I27").Value = 670.75
$ws.Range("K27").Value = 670.75
$ws.Range("M27").Value = -563.75
$ws.Range("H100").Value = 4299.6665
$ws.Range("I100").Value = 4299.6665
$ws.Range("K100").Value = 4299.6665
$ws.Range("M100").Value = -3758.6665
$ws.Range("H122").Value = 4002.4
$ws.Range("I122").Value = 3504
$ws.Range("J122").Value = 4750
$ws.Range("K122").Value = 10512
$ws.Range("L122").Value = 14250
$ws.Range("M122").Value = -8062
$ws.Range("N122").Value = -19150
$ws.Range("H132").Value = 3493.0625
$ws.Range("I132").Value = 2835.4546
$ws.Range("J132").Value = 4939.8
$ws.Range("K132").Value = 8506.363799999999
$ws.Range("L132").Value = 14819.4
$ws.Range("M132").Value = -5976.363799999999
$ws.Range("N132").Value = -19879.4
$ws.Range("H136").Value = 1600.6666
$ws.Range("I136").Value = 1600.6666
$ws.Range("K136").Value = 4801.9998
$ws.Range("M136").Value = -2251.9998

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 633.06665
$ws.Range("I100").Value = 472.81818
$ws.Range("J100").Value = 1073.75
$ws.Range("K100").Value = 945.63636
$ws.Range("L100").Value = 2147.5
$ws.Range("M100").Value = -404.63636
$ws.Range("N100").Value = -3229.5
$ws.Range("H132").Value = 1958.2667
$ws.Range("I132").Value = 1669.5714
$ws.Range("K132").Value = 5008.7142
$ws.Range("M132").Value = -2478.7142
